$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J3").Value = 2201
$ws.Range("J4").Value = 1203
$ws.Range("J8").Value = 3002
$ws.Range("J11").Value = 1119
